$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct cost for shipping AMS1117:
# - Update detailed description for the "Mach giam ap cho ESP 8266" line item
#   to reflect the added shipping cost.
$ws.Range("H7").Value = "Mach giam ap cho ESP 8266 + Shipping"

# - Update the Huy/Phu cost columns to include the shipping cost, and clear
#   the Tung column (no longer contributing to this line item).
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = $null

# Reselect the active cell as it ended up after the edit.
$ws.Range("H8").Select()
